$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price (D) and Volume (E) columns data range so that
# numeric-looking strings (e.g. "1.000", "0.9998") are preserved verbatim as text
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.574.26"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.924.20"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("D5").Value = "245.67"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "0.4885"
$ws.Range("E7").Value = "  +3.33%  "
$ws.Range("D8").Value = "0.2910"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "0.06731"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").Value = "110.67"
$ws.Range("E10").Value = "  +5.32%  "
$ws.Range("D11").Value = "19.17"
$ws.Range("E11").Value = "  +4.50%  "
$ws.Range("D12").Value = "1.921.81"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("D14").Value = "5.366"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Value = "0.6714"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "296.31"
$ws.Range("E16").Value = "  +3.19%  "
$ws.Range("D17").Value = "30.571.41"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "13.07"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").Value = "0.9998"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "5.568"
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.000007587"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").Value = "2.172.90"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "6.501"
$ws.Range("E24").Value = "  +2.85%  "
$ws.Range("D25").Value = "9.486"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("D26").Value = "164.53"
$ws.Range("E26").Value = "  -2.33%  "
$ws.Range("D27").Value = "20.26"
$ws.Range("E27").Value = "  -2.61%  "
$ws.Range("D28").Value = "2.126"
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").Value = "0.1075"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").Value = "1.450"
$ws.Range("E30").Value = "  +5.84%  "
$ws.Range("D31").Value = "4.169"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("D32").Value = "4.065"
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("D33").Value = "0.05061"
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("D34").Value = "0.7421"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  -1.06%  "
$ws.Range("D36").Value = "0.9991"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.02031"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "2.710"
$ws.Range("E38").Value = "  -1.22%  "
$ws.Range("D39").Value = "2.685"
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "110.85"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("D41").Value = "2.030"
$ws.Range("E41").Value = "  -1.51%  "
$ws.Range("D42").Value = "0.4444"
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("D43").Value = "0.8697"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("D44").Value = "71.08"
$ws.Range("E44").Value = "  +5.82%  "
$ws.Range("D45").Value = "5.841"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").Value = "0.9997"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "7.265"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").Value = "48.70"
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D49").Value = "9.204"
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("D50").Value = "0.1234"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").Value = "0.2537"
$ws.Range("E51").Value = "  +4.36%  "
